$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing article text (Science Literacy article, Hyödyllisyys column)
# was: "Käy läpi kuinka tieteellistätekstiä kannattaisi kirjoittaa 2/5"
$ws.Range("B8").Value = "Käy läpi kuinka tieteellistätekstiä kannattaisi kirjoittaa ja CA:n näkökulma asian opettamiseen 1/5"

# New article row 16: Teaching and Learning Programming and Software Engineering via Interactive Gaming
$ws.Range("A16").Value = "Teaching and Learning Programming and Software Engineering via Interactive Gaming"
$ws.Range("B16").Value = "Käsittelee MOOC sivustoa Pex4Fun (pelimäinen), pari kappaletta CA:ta piilotettuna 2/5"
$ws.Range("C16").Value = "Kerrotaan aluksi MOOC merkitys (Massive Open Online Courses). Sen jälkeen selitetään Pex4Fun sivustosta, missä voi oppia ohjelmointia. Kyseinen sivusto on pelimäinen ja samalla opettaa ohjelmointia opiskelijoille. Tehdyistä tehtävistä saa merkintöjä ylös ja voi kilpailla muitten kanssa. Duel haasteita. Pieni ripaus CA:ta miten jollakin yliopistolla professori jeesi vähän tehtävissä, jotta ei mennä liian pahasti harhaan ja masennuta."
$ws.Range("D16").Value = "Käydään läpi sivuston Pex4Fun melko selkeästi."
$ws.Range("E16").Value = "Vain muutama rivi oleellista asiaa CA:n kannalta"
$ws.Range("F16").Value = "Internet"
$ws.Range("G16").Value = "Ei tuloksia"

# New article row 18: The Abstraction Transition Taxonomy
$ws.Range("A18").Value = "The Abstraction Transition Taxonomy: Developing Desired Learning Outcomes through the Lens of Situated Cognition"
$ws.Range("B18").Value = "Käsittelee mitä ekspertiksi kehittyminen vaatii. CA ja deliperate practice. Käy läpi ohjelmointiin liittyviä taxonomeja ja kertoo niitten opetuksen tärkeydestä. Kertoo kuinka ohjelmointia opiskelevilla henkilöllä testailtiin clicker kysymyksiä (worked examples), joiden kuatta CA. 5/5"
$ws.Range("C18").Value = "Kertoo mitä ekspertiksi kehittyminen vaatii CA:n/situated cognitionin pohjalta. Kertoo lyhyesti liittyviä töitä - taxonomioiden oppiminen (sanastoa opettamiseen), Apprenticeship and Deliberate Practice (Through worked examples), Acculturation (pari programming, coding dojos). Kerrotaan oppimateriaaleista yliopiston kurssilla, jossa ~570 oppilasta. Oppimateriaalit sisältivät worked exampleja (clicker questioneja), jotka pyrkivät kysymysten kautta tuomaan esiin CA:ta - WHY kysymykset tärkeitä, how kysymysten sijaan! WHY kysymykset tuovat esiin CA:n metodeja ajatteluun ja vievät ajattelua syvemmälle (expertise), työläämpiä/hankalampia arvostella. Katsottiin myös läpi yliopiston pitämiä kurssikokeiden kysymysten muodostamisia, todella vähän why kysymyksiä."
$ws.Range("D18").Value = "Käydään worked examples näkökulmasta CA:ta. Kuinka oikein muodostetuilla kysymyksillä WHY kysymykset HOW sijaan, luodaan paremmin ammattitaitoisuutta. Kysymystä joutuu miettimään syvemmin. Uutta näkökulmaa CA:han!"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "USA, yliopisto, CS Principles course (CS0-type course), ~570 oppilasta, clicker questionit (worked examples) (english, CS speak, code)"
$ws.Range("G18").Value = "Kurssi oli onnistunut - opiskelijat kehittyivät ja näkivät teknologian mieluisampana. Abstraction Transition (AT) Taxonomy kyseisten 570 oppilaan tekemien kyselyiden perusteella (clicker questions, samalla olivat kurssimateriaaleja). Ei konkreettisia todisteita että toimiiko, mutta monet tekijät tukevat ideaa WHY kysymyksistä CA:n perusteella. "
